$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("F5").Value = 15
$ws.Range("F6").Value = 4.5

$ws.Range("F4").Select()
